$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate the last existing data row (33) three times, preserving its
# cell formatting (styles), to create rows 34-36 for the new users.
$ws.Rows.Item(33).Copy()
$ws.Rows.Item(34).Insert(-4121, 0)
$ws.Rows.Item(33).Copy()
$ws.Rows.Item(35).Insert(-4121, 0)
$ws.Rows.Item(33).Copy()
$ws.Rows.Item(36).Insert(-4121, 0)

# id
$ws.Cells.Item(34, 1).Value = 110033
$ws.Cells.Item(35, 1).Value = 110034
$ws.Cells.Item(36, 1).Value = 110035

# uin
$ws.Cells.Item(34, 2).Value = 9317596771
$ws.Cells.Item(35, 2).Value = 9317596772
$ws.Cells.Item(36, 2).Value = 9317596773

# name
$ws.Cells.Item(34, 3).Value = "Nikola Tesla"
$ws.Cells.Item(35, 3).Value = "Graham Bell"
$ws.Cells.Item(36, 3).Value = "Albert Miles"

# email
$ws.Cells.Item(34, 4).Value = "nikola.tesla@xyz.com"
$ws.Cells.Item(35, 4).Value = "graham.bell@xyz.com"
$ws.Cells.Item(36, 4).Value = "albert.miles@xyz.com"

# mobile
$ws.Cells.Item(34, 5).Value = 818876434
$ws.Cells.Item(35, 5).Value = 818876435
$ws.Cells.Item(36, 5).Value = 818876436

# Columns F (status_code), G (lang_code), H (last_login_method), I (is_active),
# J (cr_by), K (cr_dtimes) are already correct copies from row 33 (ACT / eng /
# PWD / TRUE / superadmin / now()), so no further changes are needed there.

# Reflect where the user ended up after entering the new rows - selecting the
# next empty row.
$ws.Range("A37:XFD1048576").Select() | Out-Null

Write-Output "done"
